$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 726, shifting existing rows 726:828 down to 727:829.
$ws.Rows.Item(726).Insert()

# Populate the newly inserted row 726 with the new price record.
$ws.Range("A726").Value = 10
$ws.Range("B726").Value = "Vega Modelo de Temuco"
$ws.Range("C726").Value = "La Araucanía"
$ws.Range("D726").Value = 45015
$ws.Range("E726").Value = 9
$ws.Range("F726").Value = "Fruta"
$ws.Range("G726").Value = 100102
$ws.Range("H726").Value = "Cítricos"
$ws.Range("I726").Value = 100102004
$ws.Range("J726").Value = "Mandarina"
$ws.Range("K726").Value = "Clementina"
$ws.Range("L726").Value = "Primera"
$ws.Range("M726").Value = 150
$ws.Range("N726").Value = 15000
$ws.Range("O726").Value = 15000
$ws.Range("P726").Value = 15000
$ws.Range("Q726").Value = "$/bandeja 18 kilos"
$ws.Range("R726").Value = "Región de O'Higgins"
$ws.Range("S726").Value = 833
$ws.Range("T726").Value = 18
